$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Is Active" column (F) currently stores the text string "True".
# Push the actual boolean value up instead, so the cells become real
# Excel boolean TRUE values rather than shared-string text.
$ws.Range("F2").Value = $true
$ws.Range("F3").Value = $true
